$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164; everything currently at/after row 164
# (rows 164-178) shifts down to 165-179, matching the diff exactly.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record.
$ws.Cells.Item(164, 1).Value = 2
$ws.Cells.Item(164, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(164, 3).Value = "Coquimbo"
$ws.Cells.Item(164, 4).Value = 44714
$ws.Cells.Item(164, 5).Value = 4
$ws.Cells.Item(164, 6).Value = 100112043
$ws.Cells.Item(164, 7).Value = "Pepino ensalada"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 300
$ws.Cells.Item(164, 11).Value = 17000
$ws.Cells.Item(164, 12).Value = 18000
$ws.Cells.Item(164, 13).Value = 17500
$ws.Cells.Item(164, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(164, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(164, 16).Value = 292
$ws.Cells.Item(164, 17).Value = 60
$ws.Cells.Item(164, 18).Value = "Hortaliza"
